$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the data block
# (just below the header row's accumulated history), pushing the
# previously-existing rows 174-183 down to rows 175-184.
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with the new weekly record.
$ws.Cells.Item(174, 1).Value = 7
$ws.Cells.Item(174, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(174, 3).Value = "Ñuble"
$ws.Cells.Item(174, 4).Value = 45166
$ws.Cells.Item(174, 5).Value = 16
$ws.Cells.Item(174, 6).Value = "Fruta"
$ws.Cells.Item(174, 7).Value = 100108
$ws.Cells.Item(174, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(174, 9).Value = 100108002
$ws.Cells.Item(174, 10).Value = "Mango"
$ws.Cells.Item(174, 11).Value = "Sin especificar"
$ws.Cells.Item(174, 12).Value = "Primera"
$ws.Cells.Item(174, 13).Value = 80
$ws.Cells.Item(174, 14).Value = 9000
$ws.Cells.Item(174, 15).Value = 9000
$ws.Cells.Item(174, 16).Value = 9000
$ws.Cells.Item(174, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(174, 18).Value = "Brasil"
$ws.Cells.Item(174, 19).Value = 2250
$ws.Cells.Item(174, 20).Value = 4
